# Auto-generated cell updates applying the authoritative diff of Jenova_Profits.xlsx
# Columns H..N hold computed profit-tracking figures (cached values, no live formulas).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 59719.293
$ws.Range("I28").Value = 84185.164
$ws.Range("J28").Value = 1001.2
$ws.Range("K28").Value = 84185.164
$ws.Range("L28").Value = 1001.2
$ws.Range("M28").Value = -83700.164
$ws.Range("N28").Value = -1971.2
# Row 76
$ws.Range("H76").Value = 62568280
$ws.Range("I76").Value = 171116.33
$ws.Range("K76").Value = 171116.33
$ws.Range("M76").Value = -170801.33
# Row 79
$ws.Range("H79").Value = 62568280
$ws.Range("I79").Value = 171116.33
$ws.Range("K79").Value = 171116.33
$ws.Range("M79").Value = -170024.33
# Row 138
$ws.Range("H138").Value = 6000.6895
$ws.Range("I138").Value = 3203.682
$ws.Range("J138").Value = 7709.972
$ws.Range("K138").Value = 9611.045999999998
$ws.Range("L138").Value = 23129.916
$ws.Range("M138").Value = -4471.045999999998
$ws.Range("N138").Value = -33409.916

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2411.7273
$ws.Range("J32").Value = 6497.3335
$ws.Range("L32").Value = 6497.3335
$ws.Range("N32").Value = -7071.3335
# Row 45
$ws.Range("H45").Value = 4840.8
$ws.Range("J45").Value = 8728.333000000001
$ws.Range("L45").Value = 8728.333000000001
$ws.Range("N45").Value = -9482.333000000001
# Row 61
$ws.Range("H61").Value = 6072.091
$ws.Range("I61").Value = 4965.6665
$ws.Range("K61").Value = 4965.6665
$ws.Range("M61").Value = -4753.6665
# Row 74
$ws.Range("H74").Value = 325083
$ws.Range("J74").Value = 3216.5
$ws.Range("L74").Value = 3216.5
$ws.Range("N74").Value = -4964.5
# Row 77
$ws.Range("H77").Value = 325083
$ws.Range("J77").Value = 3216.5
$ws.Range("L77").Value = 16082.5
$ws.Range("N77").Value = -24818.5
# Row 132
$ws.Range("H132").Value = 251249.2
$ws.Range("I132").Value = 451003.4
$ws.Range("K132").Value = 1353010.2
$ws.Range("M132").Value = -1350480.2
# Row 136
$ws.Range("H136").Value = 6072.091
$ws.Range("I136").Value = 4965.6665
$ws.Range("K136").Value = 14896.9995
$ws.Range("M136").Value = -12346.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 169834.33
$ws.Range("I105").Value = 1000000
$ws.Range("J105").Value = 3801.2
$ws.Range("K105").Value = 1000000
$ws.Range("L105").Value = 3801.2
$ws.Range("M105").Value = -998253
$ws.Range("N105").Value = -7295.2
# Row 107
$ws.Range("H107").Value = 1842.8572
$ws.Range("I107").Value = 1900
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1900
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 20
$ws.Range("N107").Value = -5340
# Row 134
$ws.Range("H134").Value = 96189.73
$ws.Range("I134").Value = 4298.143
$ws.Range("K134").Value = 12894.429
$ws.Range("M134").Value = -10359.429

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 711005.5600000001
$ws.Range("I31").Value = 4238713
$ws.Range("K31").Value = 4238713
$ws.Range("M31").Value = -4238418
# Row 34
$ws.Range("H34").Value = 711005.5600000001
$ws.Range("I34").Value = 4238713
$ws.Range("K34").Value = 4238713
$ws.Range("M34").Value = -4238511
# Row 50
$ws.Range("H50").Value = 47000
$ws.Range("J50").Value = 47000
$ws.Range("L50").Value = 47000
$ws.Range("N50").Value = -48250
# Row 51
$ws.Range("H51").Value = 23892.23
$ws.Range("J51").Value = 26514.143
$ws.Range("L51").Value = 26514.143
$ws.Range("N51").Value = -27986.143
# Row 60
$ws.Range("H60").Value = 29613.691
$ws.Range("J60").Value = 32543.455
$ws.Range("L60").Value = 32543.455
$ws.Range("N60").Value = -33565.455
# Row 61
$ws.Range("H61").Value = 23892.23
$ws.Range("J61").Value = 26514.143
$ws.Range("L61").Value = 26514.143
$ws.Range("N61").Value = -27210.143

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 2250.5122
$ws.Range("J68").Value = 2304.6365
$ws.Range("L68").Value = 6913.9095
$ws.Range("N68").Value = -8535.9095
# Row 71
$ws.Range("H71").Value = 2250.5122
$ws.Range("J71").Value = 2304.6365
$ws.Range("L71").Value = 20741.7285
$ws.Range("N71").Value = -28853.7285
# Row 131
$ws.Range("H131").Value = 186955.45
$ws.Range("I131").Value = 201222.4
$ws.Range("J131").Value = 175066.33
$ws.Range("K131").Value = 603667.2
$ws.Range("L131").Value = 525198.99
$ws.Range("M131").Value = -598627.2
$ws.Range("N131").Value = -535278.99

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 196.5625
$ws.Range("I2").Value = 161.37038
$ws.Range("J2").Value = 386.6
$ws.Range("K2").Value = 161.37038
$ws.Range("L2").Value = 386.6
$ws.Range("M2").Value = -48.37038000000001
$ws.Range("N2").Value = -612.6
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
# Row 97
$ws.Range("H97").Value = 1445.8096
$ws.Range("I97").Value = 1276.7333
$ws.Range("J97").Value = 1868.5
$ws.Range("K97").Value = 1276.7333
$ws.Range("L97").Value = 1868.5
$ws.Range("M97").Value = -780.7333000000001
$ws.Range("N97").Value = -2860.5
# Row 98
$ws.Range("H98").Value = 80000
$ws.Range("J98").Value = 80000
$ws.Range("L98").Value = 80000
$ws.Range("N98").Value = -85990
# Row 100
$ws.Range("H100").Value = 65000
$ws.Range("J100").Value = 65000
$ws.Range("L100").Value = 65000
$ws.Range("N100").Value = -67164
# Row 102
$ws.Range("H102").Value = 8008.7827
$ws.Range("I102").Value = 11080.167
$ws.Range("K102").Value = 11080.167
$ws.Range("M102").Value = -9458.166999999999
# Row 132
$ws.Range("H132").Value = 24659
$ws.Range("I132").Value = 2444.6177
$ws.Range("K132").Value = 7333.853099999999
$ws.Range("M132").Value = -4803.853099999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 32
$ws.Range("H32").Value = 25195
$ws.Range("I32").Value = 390
$ws.Range("J32").Value = 50000
$ws.Range("K32").Value = 390
$ws.Range("L32").Value = 50000
$ws.Range("M32").Value = -73
$ws.Range("N32").Value = -50634
# Row 68
$ws.Range("H68").Value = 70794.47
$ws.Range("I68").Value = 4136.1113
$ws.Range("K68").Value = 4136.1113
$ws.Range("M68").Value = -3387.1113
# Row 71
$ws.Range("H71").Value = 70794.47
$ws.Range("I71").Value = 4136.1113
$ws.Range("K71").Value = 20680.5565
$ws.Range("M71").Value = -16936.5565
# Row 98
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
# Row 100
$ws.Range("H100").Value = 290624.5
$ws.Range("I100").Value = 366666
$ws.Range("K100").Value = 366666
$ws.Range("M100").Value = -366125
# Row 122
$ws.Range("H122").Value = 923745.5
$ws.Range("I122").Value = 8000
$ws.Range("K122").Value = 24000
$ws.Range("M122").Value = -21550
# Row 136
$ws.Range("H136").Value = 956795.1
$ws.Range("I136").Value = 1363755
$ws.Range("K136").Value = 4091265
$ws.Range("M136").Value = -4088715

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3422.2222
$ws.Range("I81").Value = 2685.4285
$ws.Range("J81").Value = 6001
$ws.Range("K81").Value = 5370.857
$ws.Range("L81").Value = 12002
$ws.Range("M81").Value = -4309.857
$ws.Range("N81").Value = -14124
# Row 84
$ws.Range("H84").Value = 3422.2222
$ws.Range("I84").Value = 2685.4285
$ws.Range("J84").Value = 6001
$ws.Range("K84").Value = 26854.285
$ws.Range("L84").Value = 60010
$ws.Range("M84").Value = -21550.285
$ws.Range("N84").Value = -70618
# Row 100
$ws.Range("H100").Value = 1453.3529
$ws.Range("I100").Value = 1597.0714
$ws.Range("K100").Value = 3194.1428
$ws.Range("M100").Value = -2653.1428
# Row 122
$ws.Range("H122").Value = 4098.4224
$ws.Range("I122").Value = 2911.724
$ws.Range("J122").Value = 6249.3125
$ws.Range("K122").Value = 8735.172
$ws.Range("L122").Value = 18747.9375
$ws.Range("M122").Value = -6285.172
$ws.Range("N122").Value = -23647.9375
# Row 126
$ws.Range("H126").Value = 4650
$ws.Range("I126").Value = 4200
$ws.Range("K126").Value = 12600
$ws.Range("M126").Value = -10130
# Row 132
$ws.Range("H132").Value = 18218.463
$ws.Range("I132").Value = 3303.353
$ws.Range("J132").Value = 60477.945
$ws.Range("K132").Value = 9910.059000000001
$ws.Range("L132").Value = 181433.835
$ws.Range("M132").Value = -7380.059000000001
$ws.Range("N132").Value = -186493.835
# Row 133
$ws.Range("H133").Value = 65000
$ws.Range("J133").Value = 65000
$ws.Range("L133").Value = 65000
$ws.Range("N133").Value = -75120
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
